# 1.3.1.1g.xlsx - add a 2023 (column T) data point to the pensioners table,
# fix row 2's height, drop the stray Q3:S3 "no vertical align" style in favor
# of the style already used by the rest of row 3, and clear the leftover
# cell selection from the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: bump the header-separator row height, extend its border styling
#     into the new column T ---------------------------------------------
$ws.Rows.Item(2).RowHeight = 16.5

$ws.Range("T2").Value = " "
$ws.Range("T2").ClearContents()
$ws.Range("Q2").Copy()
$ws.Range("T2").PasteSpecial(-4122)

# --- Row 3: years. Q3:S3 pick up the same style already used by D3:P3
#     (vertical-centered, bordered) instead of their old un-aligned style,
#     and T3 (2023) is added with that same style ----------------------
$ws.Range("Q3").Value = 2020
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)

$ws.Range("R3").Value = 2021
$ws.Range("P3").Copy()
$ws.Range("R3").PasteSpecial(-4122)

$ws.Range("S3").Value = 2022
$ws.Range("P3").Copy()
$ws.Range("S3").PasteSpecial(-4122)

$ws.Range("T3").Value = 2023
$ws.Range("P3").Copy()
$ws.Range("T3").PasteSpecial(-4122)

# --- Row 4: number of pensioners below subsistence minimum ------------
$ws.Range("T4").Value = 263951
$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial(-4122)

# --- Row 5: share of total population ----------------------------------
$ws.Range("T5").Value = 3.7
$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial(-4122)

# --- Row 6: share of total pensioners -----------------------------------
$ws.Range("T6").Value = 32.299999999999997
$ws.Range("S6").Copy()
$ws.Range("T6").PasteSpecial(-4122)

# --- Clear the stray selection left over on the sheet view -------------
$ws.Range("A1").Select()
